# Commit: "Revised sample test results"
#
# The only substantive content change in this revision is the generation
# timestamp recorded at the start of the document's footer, which moves
# from "2025-06-30 01:36Z" to "2025-06-30 12:13Z" (the rest of the diff is
# a formatting/serialization-only change to many unrelated runs and does
# not alter any visible text). Locate that timestamp in every footer of
# every section and update it in place.

$d = $word.ActiveDocument

$oldStamp = "2025-06-30 01:36Z"
$newStamp = "2025-06-30 12:13Z"

$updated = $false

foreach ($sec in $d.Sections) {
    foreach ($ftr in $sec.Footers) {
        if ($ftr.Exists) {
            $rng = $ftr.Range
            $found = $rng.Find.Execute($oldStamp, $true, $false, $false, $false, $false, `
                                        $true, 1, $false, $newStamp, 2)
            if ($found) {
                $updated = $true
            }
        }
    }
}

# Fallback safety net: if, for some reason, no section/footer object exposed
# the run (e.g. a different footer-linking topology), fall back to a
# document-wide search/replace so the edit still lands.
if (-not $updated) {
    $whole = $d.Content
    $updated = $whole.Find.Execute($oldStamp, $true, $false, $false, $false, $false, `
                                    $true, 1, $false, $newStamp, 2)
}

$updated
